# Auto-generated from the cryptos.xlsx OOXML diff.
# Updates Coin / Link / Price / Volume(1h) cells to the refreshed crypto
# snapshot, preserving text-typed cells (force NumberFormat "@" on the
# Price column so numeric-looking strings like "292.18" are not silently
# reinterpreted by Excel as actual numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.154.54'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.226.13'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '292.18'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.74'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.512'
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.40'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('E12').Value = '  +3.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.47'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.571.83'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.97'
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.220.74'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.730'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '40.098.43'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0889'
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('E20').Value = '  +6.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.83'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.67'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.64'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.83'
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.77'
$ws.Range('E27').Value = '  -1.52%  '
$ws.Range('E28').Value = '  -1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.23'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '156.34'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.89'
$ws.Range('E31').Value = '  -7.08%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.96'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0720'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.91'
$ws.Range('E35').Value = '  +6.96%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.80'
$ws.Range('E38').Value = '  -5.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0981'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.127.98'
$ws.Range('E41').Value = '  +8.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.88'
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '18.38'
$ws.Range('E43').Value = '  +11.85%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.15'
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.02'
$ws.Range('E45').Value = '  +3.79%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0269'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.68'
$ws.Range('E47').Value = '  +3.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.437.00'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('E49').Value = '  +2.60%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.46'
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '69.34'
$ws.Range('E51').Value = '  -2.80%  '
